$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.055.23'
$ws.Range("E2").Value = '  -0.28%  '

$ws.Range("D3").Value = '1.829.16'
$ws.Range("E3").Value = '  -0.19%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9988'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6222'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.28%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.35'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.73%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07367'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.80%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2919'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.63%  '

$ws.Range("E11").Value = '  -0.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07594'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.87%  '

$ws.Range("D13").Value = '1.829.97'
$ws.Range("E13").Value = '  -0.43%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.960'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.60%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6625'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.95%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.06'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.89%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009117'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +8.99%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.024'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.29%  '

$ws.Range("D19").Value = '29.049.59'
$ws.Range("E19").Value = '  -0.35%  '

$ws.Range("D20").Value = '2.079.17'
$ws.Range("E20").Value = '  -0.42%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '225.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.94%  '

$ws.Range("E22").Value = '  -0.98%  '

$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.169'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.18%  '

$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.44'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.414'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1358'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '17.80'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.99%  '

$ws.Range("E30").Value = '  -0.79%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.055'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.47%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.035'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.23%  '

$ws.Range("E33").Value = '  +0.41%  '

$ws.Range("E34").Value = '  -1.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.840'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.58%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7342'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.35%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.151'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.20%  '

$ws.Range("E38").Value = '  +0.50%  '

$ws.Range("D39").Value = '1.283.89'
$ws.Range("E39").Value = '  +0.06%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.748'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.48%  '

$ws.Range("E41").Value = '  -0.86%  '

$ws.Range("E42").Value = '  +6.53%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8997'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.68'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.44%  '

$ws.Range("D46").Value = '1.976.99'
$ws.Range("E46").Value = '  -0.35%  '

$ws.Range("E47").Value = '  -0.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '63.85'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.90%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000120'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.44%  '

$ws.Range("E50").Value = '  -3.31%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3964'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.56%  '
